$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the unused "dimensions" column (old column G) and the blank
#     spacer row above the totals block (old row 8). Excel renumbers the
#     shared-string table and shifts formulas/refs automatically. ---
$ws.Columns.Item(7).Delete()
$ws.Rows.Item(8).Delete()

# --- New labels (added in shared-string order matching the target file) ---
$ws.Range("D9").Value = "Power budget of batteries"
$ws.Range("D10").Value = "7.2V * 180mAh * 2 = 2592mWh"
$ws.Range("I1").Value = "weight of the foam in kg"
$ws.Range("H1").Value = "volume of the foam in cm^3"
$ws.Range("F9").Value = "maximum run duration in hrs"
$ws.Range("G9").Value = "maximum run duration in minutes"

# --- New formulas: foam volume/weight budget and battery run-duration budget ---
$ws.Range("H2").Formula = "=PI() * 2.5 * (45/2)^2"
$ws.Range("I2").Formula = "=(H2/1000000)*30"
$ws.Range("F10").Formula = "=(2592/A10)"
$ws.Range("G10").Formula = "=F10*60"

# --- Column widths for the newly populated columns (bestFit-style sizing) ---
$ws.Columns.Item(6).ColumnWidth = 23.6
$ws.Columns.Item(7).ColumnWidth = 27.8
$ws.Columns.Item(8).ColumnWidth = 23.3
$ws.Columns.Item(9).ColumnWidth = 20.0

# --- Move the active selection to H2, matching the saved workbook state ---
[void]$ws.Range("H2").Select()
